$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29, pushing existing rows 29-46 down to 30-47.
# (Excel inherits the row-28 formatting for the new row, same as interactive
# "Insert Row" - this keeps column D's date number format intact.)
$ws.Rows.Item(29).Insert()

# Populate the new row 29 with the inserted record's data.
$ws.Cells.Item(29, 1).Value = 5
$ws.Cells.Item(29, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(29, 3).Value = "Maule"
$ws.Cells.Item(29, 4).Value = 44510
$ws.Cells.Item(29, 5).Value = 7
$ws.Cells.Item(29, 6).Value = 300000000
$ws.Cells.Item(29, 7).Value = "Espárragos"
$ws.Cells.Item(29, 8).Value = "Verde"
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 5000
$ws.Cells.Item(29, 11).Value = 800
$ws.Cells.Item(29, 12).Value = 800
$ws.Cells.Item(29, 13).Value = 800
$ws.Cells.Item(29, 14).Value = "`$/kilo"
$ws.Cells.Item(29, 15).Value = "Provincia de Linares"
$ws.Cells.Item(29, 16).Value = 800
$ws.Cells.Item(29, 17).Value = 1
$ws.Cells.Item(29, 18).Value = "Hortaliza"
